# ----------------------------------------------------------------------
# Apply "changes for new resolution":
#   - workbook tab-scroll ratio 993 -> 500 (view setting, best effort)
#   - add a second "frame Gen" section (two stride/kernel/padding tables)
#     below the existing table, with new blank spacer cells sprinkled
#     through the sheet matching the new look
#   - move the active selection to D59 / scroll to A21
# ----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- view / window settings (best effort; engine may not persist these) ---
try { $excel.ActiveWindow.TabRatio = 500 } catch {}
try { $excel.ActiveWindow.ScrollRow = 21 } catch {}
try { $excel.ActiveWindow.ScrollColumn = 1 } catch {}

# --- sprinkle the new blank "spacer" cells (B:F) that appear in the diff ---
# Source cell I12 already carries the plain/default style (s="0") used for
# these placeholders elsewhere in the sheet, so copy it across.
$blankSrc = $ws.Range("I12")
$blankSrc.Copy($ws.Range("B12:F13"))
$blankSrc.Copy($ws.Range("B20:F21"))
$blankSrc.Copy($ws.Range("B27:F29"))
$blankSrc.Copy($ws.Range("B35:F36"))
$blankSrc.Copy($ws.Range("B52:F53"))

# --- new section header cell ---
$ws.Range("A45").Value = "frame Gen"
$ws.Range("A45").EntireRow.RowHeight = 15

# --- first new table (rows 47-51) ---
$ws.Range("B47").Value = "in"
$ws.Range("C47").Value = "kernel"
$ws.Range("D47").Value = "padding"
$ws.Range("E47").Value = "stride"
$ws.Range("F47").Value = "out"

$ws.Range("B48").Value = 129
$ws.Range("C48").Value = 10
$ws.Range("D48").Value = 2
$ws.Range("E48").Value = 3
$ws.Range("F48").Formula = "=ROUNDDOWN((B48+2*D48-(C48-1)-1)/E48+1,0)"

$ws.Range("B49").Formula = "=F48"
$ws.Range("C49").Value = 5
$ws.Range("D49").Value = 2
$ws.Range("E49").Value = 3
$ws.Range("F49").Formula = "=ROUNDDOWN((B49+2*D49-(C49-1)-1)/E49+1,0)"

$ws.Range("B50").Formula = "=F49"
$ws.Range("C50").Value = 5
$ws.Range("D50").Value = 2
$ws.Range("E50").Value = 2
$ws.Range("F50").Formula = "=ROUNDDOWN((B50+2*D50-(C50-1)-1)/E50+1,0)"

$ws.Range("B51").Formula = "=F50"
$ws.Range("C51").Value = 5
$ws.Range("D51").Value = 2
$ws.Range("E51").Value = 2
$ws.Range("F51").Formula = "=ROUNDDOWN((B51+2*D51-(C51-1)-1)/E51+1,0)"

$ws.Range("B47:B51").EntireRow.RowHeight = 13.8

# --- second new table (rows 54-58) ---
$ws.Range("B54").Value = "in"
$ws.Range("C54").Value = "kernel"
$ws.Range("D54").Value = "padding"
$ws.Range("E54").Value = "stride"
$ws.Range("F54").Value = "out"

$ws.Range("B55").Value = 100
$ws.Range("C55").Value = 5
$ws.Range("D55").Value = 2
$ws.Range("E55").Value = 2
$ws.Range("F55").Formula = "=ROUNDDOWN((B55+2*D55-(C55-1)-1)/E55+1,0)"

$ws.Range("B56").Formula = "=F55"
$ws.Range("C56").Value = 5
$ws.Range("D56").Value = 2
$ws.Range("E56").Value = 2
$ws.Range("F56").Formula = "=ROUNDDOWN((B56+2*D56-(C56-1)-1)/E56+1,0)"

$ws.Range("B57").Formula = "=F56"
$ws.Range("C57").Value = 5
$ws.Range("D57").Value = 1
$ws.Range("E57").Value = 2
$ws.Range("F57").Formula = "=ROUNDDOWN((B57+2*D57-(C57-1)-1)/E57+1,0)"

$ws.Range("B58").Formula = "=F57"
$ws.Range("C58").Value = 5
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 2
$ws.Range("F58").Formula = "=ROUNDDOWN((B58+2*D58-(C58-1)-1)/E58+1,0)"

$ws.Range("B54:B58").EntireRow.RowHeight = 13.8

# --- update selection to match the new working area ---
$ws.Range("D59").Select()
